# Update Shadow Rate with Latest Data
#
# The shadow-rate model was refreshed against the latest vintage of data. The
# fed funds rate series (column B) is unchanged; the modeled shadow rate series
# (column C, "fedfundsrate_shadow") is recomputed for the conventional-policy
# period (rows 83-101, where it now exactly tracks the fed funds rate) and for
# the zero-lower-bound / shadow-rate period (rows 102-155), where the model
# produces updated estimates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C83").Value = 1.0099999999999998
$ws.Range("C84").Value = 1.4333333333332865
$ws.Range("C85").Value = 1.9500000000000517
$ws.Range("C86").Value = 2.46999999999995
$ws.Range("C87").Value = 2.9433333333333422
$ws.Range("C88").Value = 3.4600000000000408
$ws.Range("C89").Value = 3.980000000000028
$ws.Range("C90").Value = 4.4566666666666865
$ws.Range("C91").Value = 4.9066666666667036
$ws.Range("C92").Value = 5.2466666666666884
$ws.Range("C93").Value = 5.2466666666666884
$ws.Range("C94").Value = 5.2566666666666872
$ws.Range("C95").Value = 5.2499999999999547
$ws.Range("C96").Value = 5.0733333333333519
$ws.Range("C97").Value = 4.4966666666666821
$ws.Range("C98").Value = 3.1766666666666721
$ws.Range("C99").Value = 2.0866666666666589
$ws.Range("C100").Value = 1.9400000000000306
$ws.Range("C101").Value = 0.50666666666669968
$ws.Range("C102").Value = 1.5644100661329263
$ws.Range("C103").Value = 0.15381492810846442
$ws.Range("C104").Value = -0.75685101137021871
$ws.Range("C105").Value = -0.81554564032254584
$ws.Range("C106").Value = -0.6416259119449097
$ws.Range("C107").Value = -2.1172160245955141
$ws.Range("C108").Value = -1.9481499050362561
$ws.Range("C109").Value = -2.5944422362642161
$ws.Range("C110").Value = -2.0364267729091945
$ws.Range("C111").Value = -1.6354297542721441
$ws.Range("C112").Value = -2.9304604453149752
$ws.Range("C113").Value = -2.4986285784151963
$ws.Range("C114").Value = -3.3912490345689217
$ws.Range("C115").Value = -3.1694143390114449
$ws.Range("C116").Value = -2.7225598887246738
$ws.Range("C117").Value = -3.9315679629496891
$ws.Range("C118").Value = -2.3953712625286117
$ws.Range("C119").Value = -1.6665512287220907
$ws.Range("C120").Value = -1.2645121197149334
$ws.Range("C121").Value = -1.347818668068812
$ws.Range("C122").Value = -1.3919962197895264
$ws.Range("C123").Value = -1.3212470101792828
$ws.Range("C124").Value = -0.90436262321058614
$ws.Range("C125").Value = -0.42879133966411409
$ws.Range("C126").Value = 0.11631853045925844
$ws.Range("C127").Value = 0.011467014846999746
$ws.Range("C128").Value = 0.043484487817080364
$ws.Range("C147").Value = 7.1352820388482119
$ws.Range("C148").Value = -7.5740026381281744
$ws.Range("C149").Value = -5.1664767494255859
$ws.Range("C150").Value = -4.1916978406747134
$ws.Range("C151").Value = -3.8252569046508089
$ws.Range("C152").Value = -2.4954228601575146
$ws.Range("C153").Value = -1.3521292502786331
$ws.Range("C154").Value = -0.33951068032125331
$ws.Range("C155").Value = 0.77000000000004842
